# Add OVC_SERV.Prev.T (5-14 age band) row to the 19Tto20TMap sheet, just
# after the existing OVC_SERV.Graduated.T row (row 47), pushing every
# subsequent row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("19Tto20TMap")

# --- Insert a new row 48 (old rows 48:89 shift down to 49:90) ----------
$ws.Rows.Item(48).Insert()

# Copy the formatting of the row above (47 = OVC_SERV.Graduated.T) onto
# the new row so per-column styles (borders/fills) line up with the rest
# of the OVC_SERV block.
$ws.Range("A47:N47").Copy()
$ws.Range("A48:N48").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- Populate the new row's values --------------------------------------
$ws.Range("A48").Value = "OVC_SERV.Prev.T"
$ws.Range("B48").Value = "DE_GROUP-WTq0quAW1mf"
$ws.Range("C48").Value = "OVC_SERV"
$ws.Range("D48").Value = "RxyNwEV3oQf"
$ws.Range("E48").Value = "Numerator"
$ws.Range("F48").Value = "Som9NRMQqV7"
$ws.Range("G48").Value = "Age/Sex/ProgramStatus"
$ws.Range("H48").Value = "QG5SE83IVmL"
$ws.Range("I48").Value = "2020Oct"

# Force the age-band column to text so "5-14" isn't reinterpreted as a
# date/range by Excel's type sniffing.
$ws.Range("J48").NumberFormat = "@"
$ws.Range("J48").Value = "5-14"

$ws.Range("K48").Value = "F/M"
$ws.Range("L48").Value = "NA"
$ws.Range("M48").Value = "ag_a"
$ws.Range("N48").Value = "distribute"

# --- Refresh the autofilter range to include the new last row ----------
$ws.AutoFilterMode = $false
$ws.Range("A1:N90").AutoFilter()

# --- Keep the _FilterDatabase defined name in sync ----------------------
foreach ($n in $wb.Names) {
    if ($n.Name -eq "19Tto20TMap!_FilterDatabase") {
        $n.RefersTo = "='19Tto20TMap'!`$A`$1:`$N`$90"
    }
}

# --- Restore the active selection shown in the (scrollable) right pane -
$ws.Activate()
$ws.Range("F48").Select()
